# Updates FFXIV leve-profit market data (currentAveragePrice/LevePrice/
# LeveProfit columns) across all 8 worksheets, per scheduled market-data
# refresh. Applies scalar value updates (and clears a few now-empty cells)
# to the specific rows touched by the refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 1880.2084
$ws.Cells.Item(70, 10).Value = 2113
$ws.Cells.Item(70, 12).Value = 6339
$ws.Cells.Item(70, 14).Value = -6879

$ws.Cells.Item(73, 8).Value = 1880.2084
$ws.Cells.Item(73, 10).Value = 2113
$ws.Cells.Item(73, 12).Value = 6339
$ws.Cells.Item(73, 14).Value = -8211

$ws.Cells.Item(92, 8).Value = 366.66666
$ws.Cells.Item(92, 9).Value = 300
$ws.Cells.Item(92, 10).Value = 500
$ws.Cells.Item(92, 11).Value = 300
$ws.Cells.Item(92, 12).Value = 500
$ws.Cells.Item(92, 13).Value = 948
$ws.Cells.Item(92, 14).Value = -2996

$ws.Cells.Item(132, 8).Value = 7146.5366
$ws.Cells.Item(132, 9).Value = 3759.5938
$ws.Cells.Item(132, 11).Value = 11278.7814
$ws.Cells.Item(132, 13).Value = -8748.7814

$ws.Cells.Item(137, 8).Value = 3682.7754
$ws.Cells.Item(137, 9).Value = 4150.1577
$ws.Cells.Item(137, 10).Value = 2068.182
$ws.Cells.Item(137, 11).Value = 12450.4731
$ws.Cells.Item(137, 12).Value = 6204.545999999999
$ws.Cells.Item(137, 13).Value = -9900.473099999999
$ws.Cells.Item(137, 14).Value = -11304.546

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 669.125
$ws.Cells.Item(4, 9).Value = 411.77777
$ws.Cells.Item(4, 10).Value = 1000
$ws.Cells.Item(4, 11).Value = 411.77777
$ws.Cells.Item(4, 12).Value = 1000
$ws.Cells.Item(4, 13).Value = -295.77777
$ws.Cells.Item(4, 14).Value = -1232

$ws.Cells.Item(5, 8).Value = 334.85715
$ws.Cells.Item(5, 9).Value = 334.85715
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 334.85715
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).Value = -222.85715
$ws.Cells.Item(5, 14).ClearContents()

$ws.Cells.Item(32, 8).Value = 7186.022
$ws.Cells.Item(32, 9).Value = 5309
$ws.Cells.Item(32, 10).Value = 22577.6
$ws.Cells.Item(32, 11).Value = 5309
$ws.Cells.Item(32, 12).Value = 22577.6
$ws.Cells.Item(32, 13).Value = -5022
$ws.Cells.Item(32, 14).Value = -23151.6

$ws.Cells.Item(55, 8).Value = 24998.834
$ws.Cells.Item(55, 9).Value = 0
$ws.Cells.Item(55, 10).Value = 24998.834
$ws.Cells.Item(55, 11).Value = 0
$ws.Cells.Item(55, 12).Value = 24998.834
$ws.Cells.Item(55, 13).ClearContents()
$ws.Cells.Item(55, 14).Value = -25628.834

$ws.Cells.Item(132, 8).Value = 2598.5557
$ws.Cells.Item(132, 9).Value = 2583.1345
$ws.Cells.Item(132, 11).Value = 7749.4035
$ws.Cells.Item(132, 13).Value = -5219.4035

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(4, 8).Value = 334.85715
$ws.Cells.Item(4, 9).Value = 334.85715
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 334.85715
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = -219.85715
$ws.Cells.Item(4, 14).ClearContents()

$ws.Cells.Item(20, 8).Value = 3909.8386
$ws.Cells.Item(20, 9).Value = 4160.952
$ws.Cells.Item(20, 10).Value = 3382.5
$ws.Cells.Item(20, 11).Value = 4160.952
$ws.Cells.Item(20, 12).Value = 3382.5
$ws.Cells.Item(20, 13).Value = -3913.952
$ws.Cells.Item(20, 14).Value = -3876.5

$ws.Cells.Item(107, 8).Value = 7885.727
$ws.Cells.Item(107, 10).Value = 7561.5713
$ws.Cells.Item(107, 12).Value = 7561.5713
$ws.Cells.Item(107, 14).Value = -11401.5713

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 1862.409
$ws.Cells.Item(22, 9).Value = 1709.75
$ws.Cells.Item(22, 10).Value = 2045.6
$ws.Cells.Item(22, 11).Value = 1709.75
$ws.Cells.Item(22, 12).Value = 2045.6
$ws.Cells.Item(22, 13).Value = -1359.75
$ws.Cells.Item(22, 14).Value = -2745.6

$ws.Cells.Item(62, 8).Value = 142859790
$ws.Cells.Item(62, 9).Value = 166668880
$ws.Cells.Item(62, 11).Value = 166668880
$ws.Cells.Item(62, 13).Value = -166668256

$ws.Cells.Item(65, 8).Value = 142859790
$ws.Cells.Item(65, 9).Value = 166668880
$ws.Cells.Item(65, 11).Value = 833344400
$ws.Cells.Item(65, 13).Value = -833341280

$ws.Cells.Item(105, 8).Value = 2083.389
$ws.Cells.Item(105, 9).Value = 1896.5834
$ws.Cells.Item(105, 11).Value = 1896.5834
$ws.Cells.Item(105, 13).Value = -149.5834

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(50, 8).Value = 2201.4
$ws.Cells.Item(50, 9).Value = 365.6
$ws.Cells.Item(50, 11).Value = 1096.8
$ws.Cells.Item(50, 13).Value = -615.8000000000002

$ws.Cells.Item(53, 8).Value = 2201.4
$ws.Cells.Item(53, 9).Value = 365.6
$ws.Cells.Item(53, 11).Value = 1096.8
$ws.Cells.Item(53, 13).Value = -615.8000000000002

$ws.Cells.Item(64, 8).Value = 2934.25
$ws.Cells.Item(64, 9).Value = 2916
$ws.Cells.Item(64, 11).Value = 8748
$ws.Cells.Item(64, 13).Value = -8478

$ws.Cells.Item(67, 8).Value = 2934.25
$ws.Cells.Item(67, 9).Value = 2916
$ws.Cells.Item(67, 11).Value = 8748
$ws.Cells.Item(67, 13).Value = -7812

$ws.Cells.Item(88, 8).Value = 5746.125
$ws.Cells.Item(88, 9).Value = 3187.4
$ws.Cells.Item(88, 10).Value = 10010.667
$ws.Cells.Item(88, 11).Value = 9562.200000000001
$ws.Cells.Item(88, 12).Value = 30032.001
$ws.Cells.Item(88, 13).Value = -9134.200000000001
$ws.Cells.Item(88, 14).Value = -30888.001

$ws.Cells.Item(91, 8).Value = 5746.125
$ws.Cells.Item(91, 9).Value = 3187.4
$ws.Cells.Item(91, 10).Value = 10010.667
$ws.Cells.Item(91, 11).Value = 9562.200000000001
$ws.Cells.Item(91, 12).Value = 30032.001
$ws.Cells.Item(91, 13).Value = -8080.200000000001
$ws.Cells.Item(91, 14).Value = -32996.001

$ws.Cells.Item(98, 8).Value = 1667271.6
$ws.Cells.Item(98, 10).Value = 655
$ws.Cells.Item(98, 12).Value = 1965
$ws.Cells.Item(98, 14).Value = -4961

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 56.5625
$ws.Cells.Item(2, 9).Value = 35.083332
$ws.Cells.Item(2, 11).Value = 35.083332
$ws.Cells.Item(2, 13).Value = 77.916668

$ws.Cells.Item(62, 8).Value = 99989.5
$ws.Cells.Item(62, 10).Value = 99989.5
$ws.Cells.Item(62, 12).Value = 99989.5
$ws.Cells.Item(62, 14).Value = -101361.5

$ws.Cells.Item(65, 8).Value = 99989.5
$ws.Cells.Item(65, 10).Value = 99989.5
$ws.Cells.Item(65, 12).Value = 299968.5
$ws.Cells.Item(65, 14).Value = -306832.5

$ws.Cells.Item(113, 8).Value = 3811
$ws.Cells.Item(113, 9).Value = 3811
$ws.Cells.Item(113, 11).Value = 3811
$ws.Cells.Item(113, 13).Value = -1641

$ws.Cells.Item(126, 8).Value = 5638.3335
$ws.Cells.Item(126, 9).Value = 6401.1665
$ws.Cells.Item(126, 11).Value = 19203.4995
$ws.Cells.Item(126, 13).Value = -16733.4995

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2122.2
$ws.Cells.Item(7, 9).Value = 2122.2
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 2122.2
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 13).Value = -2010.2
$ws.Cells.Item(7, 14).ClearContents()

$ws.Cells.Item(16, 8).Value = 8337376.5
$ws.Cells.Item(16, 9).Value = 8932189
$ws.Cells.Item(16, 11).Value = 8932189
$ws.Cells.Item(16, 13).Value = -8932019

$ws.Cells.Item(22, 8).Value = 952.3
$ws.Cells.Item(22, 9).Value = 961.5714
$ws.Cells.Item(22, 11).Value = 961.5714
$ws.Cells.Item(22, 13).Value = -666.5714

$ws.Cells.Item(27, 8).Value = 952.3
$ws.Cells.Item(27, 9).Value = 961.5714
$ws.Cells.Item(27, 11).Value = 961.5714
$ws.Cells.Item(27, 13).Value = -854.5714

$ws.Cells.Item(40, 8).Value = 0
$ws.Cells.Item(40, 9).Value = 0
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 0
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 13).ClearContents()
$ws.Cells.Item(40, 14).ClearContents()

$ws.Cells.Item(46, 8).Value = 2974.577
$ws.Cells.Item(46, 9).Value = 2659
$ws.Cells.Item(46, 10).Value = 3206
$ws.Cells.Item(46, 11).Value = 2659
$ws.Cells.Item(46, 12).Value = 3206
$ws.Cells.Item(46, 13).Value = -2471
$ws.Cells.Item(46, 14).Value = -3582

$ws.Cells.Item(55, 8).Value = 1654.4193
$ws.Cells.Item(55, 9).Value = 1471
$ws.Cells.Item(55, 10).Value = 2039.6
$ws.Cells.Item(55, 11).Value = 1471
$ws.Cells.Item(55, 12).Value = 2039.6
$ws.Cells.Item(55, 13).Value = -1298
$ws.Cells.Item(55, 14).Value = -2385.6

$ws.Cells.Item(126, 8).Value = 2122.2
$ws.Cells.Item(126, 9).Value = 2122.2
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 6366.599999999999
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).Value = -3896.599999999999
$ws.Cells.Item(126, 14).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 3994.4167
$ws.Cells.Item(4, 10).Value = 5304.25
$ws.Cells.Item(4, 12).Value = 5304.25
$ws.Cells.Item(4, 14).Value = -5530.25

$ws.Cells.Item(113, 8).Value = 689.73334
$ws.Cells.Item(113, 9).Value = 394
$ws.Cells.Item(113, 11).Value = 1182
$ws.Cells.Item(113, 13).Value = 988
